$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room for the new article at row 2 by shifting all the existing data rows
#    down by one (copying values from the bottom up). Using plain value copies instead
#    of Rows.Insert() avoids Excel auto-copying the bold header style onto row 2.
for ($r = 21; $r -ge 2; $r--) {
    $dst = $r + 1
    $ws.Range("A$dst").Value = $ws.Range("A$r").Value2
    $ws.Range("B$dst").Value = $ws.Range("B$r").Value2
    $ws.Range("C$dst").Value = $ws.Range("C$r").Value2
}

# 2) Populate the brand-new article row (row 2), overwriting the row that still
#    holds a duplicate of the old row 2 content after the shift above.
$ws.Range("A2").Value = 'Hà Nội cho học sinh nghỉ học ngày 6/10'
$ws.Range("B2").Value = '2025-10-05T18:24:00'
$ws.Range("C2").Value = 'Sở GD&ĐT Hà Nội đề nghị các đơn vị cho học sinh nghỉ học ngày mai (6/10) để tránh bão số 11 Matmo.'
$ws.Range("D2").Value = 'https://cafef.vn/ha-noi-cho-23-trieu-hoc-sinh-nghi-hoc-ngay-6-10-ung-pho-bao-so-11-matmo-188251005182149527.chn'

# 3) Add the "Link" column value (column D) to every pre-existing article row,
#    which have now shifted down to rows 3-21.
$ws.Range("D3").Value = 'https://cafef.vn/khong-co-viec-hoang-huong-van-dang-lam-tu-thien-tren-facebook-du-da-bi-khoi-to-bat-tam-giam-188251005165823089.chn'
$ws.Range("D4").Value = 'https://cafef.vn/ong-le-ngoc-quang-duoc-chi-dinh-giu-chuc-bi-thu-tinh-uy-quang-tri-188251005165602554.chn'
$ws.Range("D5").Value = 'https://cafef.vn/bat-tam-giam-thuy-tet-be-nho-188251005142828335.chn'
$ws.Range("D6").Value = 'https://cafef.vn/kham-xet-va-bat-tam-giam-giam-doc-bui-van-binh-188251005131259582.chn'
$ws.Range("D7").Value = 'https://cafef.vn/thong-tin-khan-ve-5-san-bay-anh-huong-bao-matmo-188251005131457805.chn'
$ws.Range("D8").Value = 'https://cafef.vn/ben-trong-xuong-san-xuat-22-trieu-binh-chua-chay-gia-cong-thuc-pha-tron-tu-nghi-ra-54-la-bot-khong-co-tac-dung-dap-lua-188251005113331212.chn'
$ws.Range("D9").Value = 'https://cafef.vn/vu-gay-roi-tai-karaoke-dubai-o-da-nang-bat-xin-nghiem-cung-3-nguoi-18825100511362671.chn'
$ws.Range("D10").Value = 'https://cafef.vn/chu-tich-xa-noi-gi-ve-viec-bi-thu-xa-bi-to-bo-lop-boi-duong-de-choi-pickleball-188251005113845969.chn'
$ws.Range("D11").Value = 'https://cafef.vn/quang-ninh-cam-bien-tu-8h-ngay-5-10-188251005103156113.chn'
$ws.Range("D12").Value = 'https://cafef.vn/bao-so-11-tang-cap-di-nhu-luot-tren-bien-khi-nao-se-do-bo-anh-huong-nhung-vung-nao-cua-viet-nam-188251005094305869.chn'
$ws.Range("D13").Value = 'https://cafef.vn/chu-tich-ha-noi-chu-dong-cho-hoc-sinh-nghi-hoc-nguoi-lao-dong-lam-viec-online-ung-pho-bao-so-11-18825100510175321.chn'
$ws.Range("D14").Value = 'https://cafef.vn/cong-an-bat-qua-tang-hoang-van-tuyen-khi-dang-livestream-188251005094055791.chn'
$ws.Range("D15").Value = 'https://cafef.vn/cong-an-canh-bao-khan-ve-thong-tin-chinh-phu-tang-moi-nguoi-dan-1-trieu-dong-dip-tet-188251005083649401.chn'
$ws.Range("D16").Value = 'https://cafef.vn/ha-noi-khuyen-khich-nguoi-dan-lam-viec-online-vao-thu-hai-ngay-6-10-188251005084453025.chn'
$ws.Range("D17").Value = 'https://cafef.vn/thi-hanh-lenh-kham-xet-noi-o-cua-nu-ke-toan-pham-thi-thuy-188251005083240511.chn'
$ws.Range("D18").Value = 'https://cafef.vn/bat-tam-giam-giam-doc-phan-the-hoai-va-em-trai-phan-bach-thong-188251005071021839.chn'
$ws.Range("D19").Value = 'https://cafef.vn/bao-so-11-suc-gio-rat-manh-sap-do-bo-vao-nuoc-ta-ha-noi-va-16-tinh-thanh-dac-biet-chu-y-188251005071318815.chn'
$ws.Range("D20").Value = 'https://cafef.vn/bat-tam-giam-kham-xet-noi-o-cua-nguyen-van-trinh-188251005070837688.chn'
$ws.Range("D21").Value = 'https://cafef.vn/ha-noi-chay-phong-giao-dich-ngan-hang-o-duong-thuy-khue-188251004213035091.chn'

# 4) Remove the now-last, oldest article row (previously row 21, shifted to row 22)
$ws.Rows.Item(22).Delete()

# 5) Add the "Link" header in D1, matching the bold/bordered header style used by A1:C1
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Link"

# 6) Size the new column to fit its contents
$ws.Columns.Item(4).ColumnWidth = 140.2

# 7) Update the sheet view: drop the old frozen/scrolled topLeftCell and update the selection
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A15:A24").Select()
